# The last paragraph of the document (right before the sectPr) is an
# empty, numbered "ListParagraph" item that only carries the
# now-orphaned "_GoBack" bookmark:
#
#   <w:p>
#     <w:pPr>
#       <w:pStyle w:val="ListParagraph"/>
#       <w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr>
#       <w:rPr>...Times New Roman, 13pt...</w:rPr>
#     </w:pPr>
#     <w:bookmarkStart w:id="0" w:name="_GoBack"/>
#     <w:bookmarkEnd w:id="0"/>
#   </w:p>
#
# It gets replaced by seven plain (non-list) paragraphs: two blank
# spacer paragraphs, then five paragraphs of new "Professional
# contribution" notes (the bookmark now lives in the third one, next to
# the "What the thesis do:" label).

$d = $word.ActiveDocument

$rPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr>'
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function PlainPara([string]$text) {
    return "<w:p $wNs><w:pPr>$rPr</w:pPr><w:r>$rPr<w:t>$text</w:t></w:r></w:p>"
}

$blank1 = "<w:p $wNs><w:pPr>$rPr</w:pPr></w:p>"
$blank2 = "<w:p $wNs><w:pPr>$rPr</w:pPr></w:p>"

$withBookmark = "<w:p $wNs><w:pPr>$rPr</w:pPr><w:r>$rPr<w:t>What the thesis do:</w:t></w:r>" +
                '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$dataset  = PlainPara "Build a train dataset. Then, recognizing emotion."
$annSvm   = PlainPara "ANN or SVM?"
$asmAam   = PlainPara "ASM or AAM?"
$compare  = PlainPara "Run each technique =&gt; Compare the results =&gt; Choose one fitss"

$newXml = $blank1 + $blank2 + $withBookmark + $dataset + $annSvm + $asmAam + $compare

$target = $d.Paragraphs.Last
$target.Range.InsertXML($newXml)
